# Add new 5G propagation model tests
#
# Updates the measured path_loss (dB) inputs for the 4G baseline scenario
# (row 3) and the three 5G scenarios (rows 9, 12, 15) on the
# "link_buget_example 1" sheet. Every downstream formula (SINR, capacity,
# energy-efficiency calculations further down the sheet) recalculates
# automatically from these inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("link_buget_example 1")
$ws.Activate()

# path_loss (dB) input cells (column E) for the four link-budget scenarios
$ws.Range("E3").Value = 111
$ws.Range("E9").Value = 124
$ws.Range("E12").Value = 113
$ws.Range("E15").Value = 129

# Leave the sheet scrolled to the top with B21 selected, matching the
# author's final view of the workbook after adding the new test cases.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("B21").Select() | Out-Null
